$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Round Tanks" worksheet, placed after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Round Tanks"

# --- Text labels, entered in the same order the author originally typed them
#     (this keeps the shared-string table ordering faithful to the source) ---
$ws2.Range("B1").Value = "Diameter"
$ws2.Range("B3").Value = "Height"
$ws2.Range("B11").Value = "(multiplier)"
$ws2.Range("A13").Value = "Entry cost"
$ws2.Range("A14").Value = "cost"
$ws2.Range("A15").Value = "mass"
$ws2.Range("B4").Value = "Area"
$ws2.Range("B2").Value = "Radius"
$ws2.Range("B5").Value = "Volume"
$ws2.Range("B6").Value = "Total Surface Area"
$ws2.Range("A16").Value = "K&B"
$ws2.Range("B7").Value = "K&B/unit volume"

# --- Row 1 - Diameter (input values, inches) ---
$ws2.Range("D1").Value = 1.25
$ws2.Range("E1").Value = 1.875
$ws2.Range("F1").Value = 2.5
$ws2.Range("G1").Value = 3.75
$ws2.Range("H1").Value = 5

# --- Row 2 - Radius ---
$ws2.Range("D2").Formula = '=D1/2'
$ws2.Range("E2").Formula = '=E1/2'
$ws2.Range("F2").Formula = '=F1/2'
$ws2.Range("G2").Formula = '=G1/2'
$ws2.Range("H2").Formula = '=H1/2'

# --- Row 3 - Height ---
$ws2.Range("D3").Formula = '=D1/3'
$ws2.Range("E3").Formula = '=E11*$D$3'
$ws2.Range("F3").Formula = '=F11*$D$3'
$ws2.Range("G3").Formula = '=G11*$D$3'
$ws2.Range("H3").Formula = '=H11*$D$3'

# --- Row 4 - Area ---
$ws2.Range("D4").Formula = '=PI()*D2^2'
$ws2.Range("E4").Formula = '=PI()*E2^2'
$ws2.Range("F4").Formula = '=PI()*F2^2'
$ws2.Range("G4").Formula = '=PI()*G2^2'
$ws2.Range("H4").Formula = '=PI()*H2^2'

# --- Row 5 - Volume ---
$ws2.Range("D5").Formula = '=PI()*D2^2*D3'
$ws2.Range("E5").Formula = '=PI()*E2^2*E3'
$ws2.Range("F5").Formula = '=PI()*F2^2*F3'
$ws2.Range("G5").Formula = '=PI()*G2^2*G3'
$ws2.Range("H5").Formula = '=PI()*H2^2*H3'

# --- Row 6 - Total Surface Area ---
$ws2.Range("D6").Formula = '=2*PI()*(D1/2)*D3+2*PI()*D2^2'
$ws2.Range("E6").Formula = '=2*PI()*(E1/2)*E3+2*PI()*E2^2'
$ws2.Range("F6").Formula = '=2*PI()*(F1/2)*F3+2*PI()*F2^2'
$ws2.Range("G6").Formula = '=2*PI()*(G1/2)*G3+2*PI()*G2^2'
$ws2.Range("H6").Formula = '=2*PI()*(H1/2)*H3+2*PI()*H2^2'

# --- Row 7 - K&B/unit volume ---
$ws2.Range("D7").Formula = '=D16/D5'

# --- Row 11 - (multiplier) ---
$ws2.Range("E11").Formula = '=E1/$D$1'
$ws2.Range("F11").Formula = '=F1/$D$1'
$ws2.Range("G11").Formula = '=G1/$D$1'
$ws2.Range("H11").Formula = '=H1/$D$1'

# --- Row 13 - Entry cost ---
$ws2.Range("D13").Value = 2500
$ws2.Range("E13").Formula = '=E11^2*$D$14'
$ws2.Range("F13").Formula = '=F11^2*$D$14'
$ws2.Range("G13").Formula = '=G11^2*$D$14'
$ws2.Range("H13").Formula = '=H11^2*$D$14'

# --- Row 14 - cost ---
$ws2.Range("D14").Value = 1200
$ws2.Range("E14").Formula = '=$E$11*D14'
$ws2.Range("F14").Formula = '=$E$11*E14'
$ws2.Range("G14").Formula = '=$E$11*F14'
$ws2.Range("H14").Formula = '=$E$11*G14'

# --- Row 15 - mass ---
$ws2.Range("D15").Value = 0.06
$ws2.Range("E15").Formula = '=$D$15/$D$6*E6'
$ws2.Range("F15").Formula = '=$D$15/$D$6*F6'
$ws2.Range("G15").Formula = '=$D$15/$D$6*G6'
$ws2.Range("H15").Formula = '=$D$15/$D$6*H6'

# --- Row 16 - K&B ---
$ws2.Range("D16").Value = 50
$ws2.Range("E16").Formula = '=E5*$D$7'
$ws2.Range("F16").Formula = '=F5*$D$7'
$ws2.Range("G16").Formula = '=G5*$D$7'
$ws2.Range("H16").Formula = '=H5*$D$7'

# Match the portrait page orientation used throughout the workbook
$ws2.PageSetup.Orientation = 1

# Selection / view state for the new sheet, and make it the active tab
$ws2.Range("H16").Select()
$ws2.Activate()
